# update database and change read_price algorithm
#
# The quarterly report rolls its 10-quarter trailing window forward by one
# quarter: the oldest quarter (1399/06) drops off the front and a new
# quarter (1401/12) is appended at the end. Every per-quarter figure in
# columns E:N shifts one column to the left and a freshly reported value is
# written into column N.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Column headers (row 8 and row 24 share the same 10 quarter labels) ---
$quarterLabels = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

for ($i = 0; $i -lt $quarterLabels.Length; $i++) {
    $col = 5 + $i   # E=5 .. N=14
    $ws.Cells.Item(8, $col).Value = $quarterLabels[$i]
    $ws.Cells.Item(24, $col).Value = $quarterLabels[$i]
}

# --- Data rows: shift left by one quarter, append the newly reported value ---
$rowData = @{
    10 = @(54080, 67558, 47833, 76948, 71523, 82679, 65471, 88111, 170011, 244157)
    13 = @(30852, 71767, 109256, 94441, 327662, -45914, 125403, 288427, 396937, 271104)
    16 = @(2949, 11298, 4286, 3975, 4131, 4356, 6045, 2659, 1089, 18657)
    17 = @(259763, 336000, 414808, 440661, 485089, 558406, 581398, 742007, 659850, 462170)
    19 = @(85859, 128677, 81457, 142317, 14470, 166243, 162842, 47449, 109194, 322341)
    20 = @(433503, 615300, 657640, 758342, 902875, 765770, 941159, 1168653, 1337081, 1318429)
    26 = @(1155, 1190, 1208, 1263, 1245, 1235, 1209, 1250, 1215, 1300)
    27 = @(556, 578, 591, 618, 581, 571, 548, 555, 552, 600)
}

foreach ($row in $rowData.Keys) {
    $values = $rowData[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 5 + $i   # E=5 .. N=14
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
